# Update the team-specific transition-probability matrix on Sheet1.
# These cells hold empirical transition probabilities (each data row sums to 1);
# additional simulated games shifted the counts, so the probabilities below
# reflect the refreshed game totals referenced in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 0.07692307692307693
$ws.Range("C2").Value = 0.6923076923076923
$ws.Range("J2").Value = 0.1538461538461539
$ws.Range("P2").Value = 0.07692307692307693

# Row 3
$ws.Range("J3").Value = 0.1111111111111111
$ws.Range("P3").Value = 0.5555555555555556
$ws.Range("S3").Value = 0.3333333333333333

# Row 4
$ws.Range("J4").Value = 0.1111111111111111
$ws.Range("P4").Value = 0.3333333333333333
$ws.Range("S4").Value = 0.5555555555555556

# Row 6
$ws.Range("Q6").Value = 0.5555555555555556
$ws.Range("S6").Value = 0.4444444444444444

# Row 7
$ws.Range("B7").Value = 0.2
$ws.Range("O7").Value = 0.2
$ws.Range("Q7").Value = 0.2
$ws.Range("S7").Value = 0.4

# Row 8
$ws.Range("B8").Value = 0.04166666666666666
$ws.Range("E8").Value = 0.04166666666666666
$ws.Range("F8").Value = 0.08333333333333333
$ws.Range("J8").Value = 0.125
$ws.Range("R8").Value = 0.125
$ws.Range("S8").Value = 0.25

# Row 9
$ws.Range("B9").Value = 0.1052631578947368
$ws.Range("D9").Value = 0.05263157894736842
$ws.Range("F9").Value = 0.1052631578947368
$ws.Range("J9").Value = 0.05263157894736842
$ws.Range("Q9").Value = 0.3684210526315789
$ws.Range("S9").Value = 0.3157894736842105

# Row 10
$ws.Range("B10").Value = 0.07079646017699115
$ws.Range("D10").Value = 0.07079646017699115
$ws.Range("F10").Value = 0.02654867256637168
$ws.Range("J10").Value = 0.1769911504424779
$ws.Range("O10").Value = 0.01769911504424779
$ws.Range("Q10").Value = 0.2743362831858407
$ws.Range("R10").Value = 0.09734513274336283
$ws.Range("S10").Value = 0.2654867256637168

# Row 11
$ws.Range("G11").Value = 0.2
$ws.Range("K11").Value = 0.2
$ws.Range("L11").Value = 0.6

# Row 13
$ws.Range("G13").Value = 1

# Row 15
$ws.Range("H15").Value = 0.1176470588235294
$ws.Range("I15").Value = 0.05882352941176471
$ws.Range("J15").Value = 0.5294117647058824
$ws.Range("O15").Value = 0.05882352941176471
$ws.Range("S15").Value = 0.2352941176470588

# Row 16
$ws.Range("H16").Value = 0.1
$ws.Range("J16").Value = 0.8
$ws.Range("S16").Value = 0.1

# Row 17
$ws.Range("H17").Value = 0.1153846153846154
$ws.Range("I17").Value = 0.1538461538461539
$ws.Range("J17").Value = 0.4230769230769231
$ws.Range("K17").Value = 0.03846153846153846
$ws.Range("O17").Value = 0.07692307692307693
$ws.Range("S17").Value = 0.1923076923076923

# Row 18
$ws.Range("H18").Value = 0.1428571428571428
$ws.Range("I18").Value = 0.1428571428571428
$ws.Range("J18").Value = 0.5714285714285714
$ws.Range("O18").Value = 0.1428571428571428

# Row 19
$ws.Range("H19").Value = 0.1794871794871795
$ws.Range("I19").Value = 0.1025641025641026
$ws.Range("J19").Value = 0.5128205128205128
$ws.Range("K19").Value = 0.02564102564102564
$ws.Range("M19").Value = 0.01282051282051282
$ws.Range("O19").Value = 0.0641025641025641
$ws.Range("S19").Value = 0.1025641025641026
